# Week 10 grade update: fill in the previously-blank "work/comment" pair
# (columns T/U) for four students, and move the active-cell selection.
#
# Columns T (grade, e.g. "5.0"/"4.5") hold values that *look* numeric.
# A plain `.Value = "5.0"` assignment gets silently coerced to the number
# 5 by Excel's type inference, which would store it as a numeric cell
# instead of the shared-string text cell the workbook actually uses
# elsewhere (see D3, F3, ... all t="s"). To force a genuine text value
# without perturbing the cell's style (setting NumberFormat to force text
# leaves an orphan style entry behind), each cell is first given a
# literal-string formula, then "frozen" to a static value via a
# values-only self-paste, and finally re-stamped with the style of a
# sibling cell (T4) that already carries the correct "has a grade" style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-GradeText {
    # Positional params: this PS host doesn't bind named (-Foo) args
    # reliably inside functions, so call this with positional args only.
    param($Cell, $Text, $StyleDonor)
    # Compute the literal text via a formula, then collapse the formula
    # down to a plain cached value (still text, no auto-number coercion).
    $ws.Range($Cell).Formula = '="' + $Text + '"'
    $ws.Range($Cell).Copy()
    $ws.Range($Cell).PasteSpecial(-4163)  # xlPasteValues

    # Re-apply the "filled in" style (matches other graded T cells)
    # without touching the value we just wrote.
    $ws.Range($StyleDonor).Copy()
    $ws.Range($Cell).PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = 0

# Row 5 - Samuel Castañeda Montoya
Set-GradeText "T5" "5.0" "T4"
$ws.Range("U5").Value = "#2"

# Row 7 - Jose Fernando Maya Ramirez
Set-GradeText "T7" "4.5" "T4"
$ws.Range("U7").Value = "#1"

# Row 10 - Sebastian Rendon Rendon
Set-GradeText "T10" "5.0" "T4"
$ws.Range("U10").Value = "#1"

# Row 14 - Mateo Velasquez Zapata
Set-GradeText "T14" "5.0" "T4"
$ws.Range("U14").Value = "#3"

$excel.CutCopyMode = 0

# Move the active cell / selection (was U8 -> now U1)
$ws.Range("U1").Select()
